$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 213.2
$ws.Cells.Item(12, 9).Value = 215.55556
$ws.Cells.Item(12, 10).Value = 192
$ws.Cells.Item(12, 11).Value = 215.55556
$ws.Cells.Item(12, 12).Value = 192
$ws.Cells.Item(12, 13).Value = -45.55556000000001
$ws.Cells.Item(12, 14).Value = -532

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3496.0881
$ws.Cells.Item(64, 9).Value = 3485.1333
$ws.Cells.Item(64, 10).Value = 3504.7368
$ws.Cells.Item(64, 11).Value = 3485.1333
$ws.Cells.Item(64, 12).Value = 3504.7368
$ws.Cells.Item(64, 13).Value = -3237.1333
$ws.Cells.Item(64, 14).Value = -4000.7368

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 3496.0881
$ws.Cells.Item(67, 9).Value = 3485.1333
$ws.Cells.Item(67, 10).Value = 3504.7368
$ws.Cells.Item(67, 11).Value = 3485.1333
$ws.Cells.Item(67, 12).Value = 3504.7368
$ws.Cells.Item(67, 13).Value = -2627.1333
$ws.Cells.Item(67, 14).Value = -5220.736800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 73921.766
$ws.Cells.Item(74, 9).Value = 228745.75
$ws.Cells.Item(74, 10).Value = 5111.1113
$ws.Cells.Item(74, 11).Value = 228745.75
$ws.Cells.Item(74, 12).Value = 5111.1113
$ws.Cells.Item(74, 13).Value = -227809.75
$ws.Cells.Item(74, 14).Value = -6983.1113

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 5369.143
$ws.Cells.Item(76, 9).Value = 6437.875
$ws.Cells.Item(76, 10).Value = 3944.1667
$ws.Cells.Item(76, 11).Value = 6437.875
$ws.Cells.Item(76, 12).Value = 3944.1667
$ws.Cells.Item(76, 13).Value = -6122.875
$ws.Cells.Item(76, 14).Value = -4574.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(77, 8).Value = 73921.766
$ws.Cells.Item(77, 9).Value = 228745.75
$ws.Cells.Item(77, 10).Value = 5111.1113
$ws.Cells.Item(77, 11).Value = 1143728.75
$ws.Cells.Item(77, 12).Value = 25555.5565
$ws.Cells.Item(77, 13).Value = -1139048.75
$ws.Cells.Item(77, 14).Value = -34915.5565

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(79, 8).Value = 5369.143
$ws.Cells.Item(79, 9).Value = 6437.875
$ws.Cells.Item(79, 10).Value = 3944.1667
$ws.Cells.Item(79, 11).Value = 6437.875
$ws.Cells.Item(79, 12).Value = 3944.1667
$ws.Cells.Item(79, 13).Value = -5345.875
$ws.Cells.Item(79, 14).Value = -6128.1667

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 198.42857
$ws.Cells.Item(5, 9).Value = 214.83333
$ws.Cells.Item(5, 10).Value = 100
$ws.Cells.Item(5, 11).Value = 214.83333
$ws.Cells.Item(5, 12).Value = 100
$ws.Cells.Item(5, 13).Value = -102.83333
$ws.Cells.Item(5, 14).Value = -324

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1842.5111
$ws.Cells.Item(61, 9).Value = 1740.3422
$ws.Cells.Item(61, 10).Value = 2397.1428
$ws.Cells.Item(61, 11).Value = 1740.3422
$ws.Cells.Item(61, 12).Value = 2397.1428
$ws.Cells.Item(61, 13).Value = -1528.3422
$ws.Cells.Item(61, 14).Value = -2821.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 3301
$ws.Cells.Item(63, 9).Value = 2626.25
$ws.Cells.Item(63, 11).Value = 2626.25
$ws.Cells.Item(63, 13).Value = -1940.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 3301
$ws.Cells.Item(66, 9).Value = 2626.25
$ws.Cells.Item(66, 11).Value = 13131.25
$ws.Cells.Item(66, 13).Value = -9699.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 1842.5111
$ws.Cells.Item(136, 9).Value = 1740.3422
$ws.Cells.Item(136, 10).Value = 2397.1428
$ws.Cells.Item(136, 11).Value = 5221.0266
$ws.Cells.Item(136, 12).Value = 7191.428400000001
$ws.Cells.Item(136, 13).Value = -2671.0266
$ws.Cells.Item(136, 14).Value = -12291.4284

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 198.42857
$ws.Cells.Item(4, 9).Value = 214.83333
$ws.Cells.Item(4, 10).Value = 100
$ws.Cells.Item(4, 11).Value = 214.83333
$ws.Cells.Item(4, 12).Value = 100
$ws.Cells.Item(4, 13).Value = -99.83332999999999
$ws.Cells.Item(4, 14).Value = -330

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2221.0605
$ws.Cells.Item(134, 9).Value = 1677
$ws.Cells.Item(134, 10).Value = 3309.182
$ws.Cells.Item(134, 11).Value = 5031
$ws.Cells.Item(134, 12).Value = 9927.545999999998
$ws.Cells.Item(134, 13).Value = -2496
$ws.Cells.Item(134, 14).Value = -14997.546

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1631.0869
$ws.Cells.Item(31, 9).Value = 1338.6389
$ws.Cells.Item(31, 11).Value = 1338.6389
$ws.Cells.Item(31, 13).Value = -1043.6389

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 1631.0869
$ws.Cells.Item(34, 9).Value = 1338.6389
$ws.Cells.Item(34, 11).Value = 1338.6389
$ws.Cells.Item(34, 13).Value = -1136.6389

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 863099.2
$ws.Cells.Item(58, 9).Value = 1483129.1
$ws.Cells.Item(58, 10).Value = 1946.4445
$ws.Cells.Item(58, 11).Value = 1483129.1
$ws.Cells.Item(58, 12).Value = 1946.4445
$ws.Cells.Item(58, 13).Value = -1482926.1
$ws.Cells.Item(58, 14).Value = -2352.4445

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 102761
$ws.Cells.Item(62, 9).Value = 168601.67
$ws.Cells.Item(62, 10).Value = 4000
$ws.Cells.Item(62, 11).Value = 168601.67
$ws.Cells.Item(62, 12).Value = 4000
$ws.Cells.Item(62, 13).Value = -167977.67
$ws.Cells.Item(62, 14).Value = -5248

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 102761
$ws.Cells.Item(65, 9).Value = 168601.67
$ws.Cells.Item(65, 10).Value = 4000
$ws.Cells.Item(65, 11).Value = 843008.3500000001
$ws.Cells.Item(65, 12).Value = 20000
$ws.Cells.Item(65, 13).Value = -839888.3500000001
$ws.Cells.Item(65, 14).Value = -26240

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 863099.2
$ws.Cells.Item(136, 9).Value = 1483129.1
$ws.Cells.Item(136, 10).Value = 1946.4445
$ws.Cells.Item(136, 11).Value = 4449387.300000001
$ws.Cells.Item(136, 12).Value = 5839.333500000001
$ws.Cells.Item(136, 13).Value = -4446837.300000001
$ws.Cells.Item(136, 14).Value = -10939.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 634.7568
$ws.Cells.Item(113, 10).Value = 680.375
$ws.Cells.Item(113, 12).Value = 2041.125
$ws.Cells.Item(113, 14).Value = -6381.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(118, 8).Value = 2236.9285
$ws.Cells.Item(118, 9).Value = 707.25
$ws.Cells.Item(118, 10).Value = 2848.8
$ws.Cells.Item(118, 11).Value = 2121.75
$ws.Cells.Item(118, 12).Value = 8546.400000000001
$ws.Cells.Item(118, 13).Value = -878.75
$ws.Cells.Item(118, 14).Value = -11032.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 9).Value = 456
$ws.Cells.Item(131, 11).Value = 1368
$ws.Cells.Item(131, 13).Value = 3672

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5513.136
$ws.Cells.Item(70, 9).Value = 5416.0557
$ws.Cells.Item(70, 10).Value = 5950
$ws.Cells.Item(70, 11).Value = 5416.0557
$ws.Cells.Item(70, 12).Value = 5950
$ws.Cells.Item(70, 13).Value = -5146.0557
$ws.Cells.Item(70, 14).Value = -6490

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 5513.136
$ws.Cells.Item(73, 9).Value = 5416.0557
$ws.Cells.Item(73, 10).Value = 5950
$ws.Cells.Item(73, 11).Value = 5416.0557
$ws.Cells.Item(73, 12).Value = 5950
$ws.Cells.Item(73, 13).Value = -4480.0557
$ws.Cells.Item(73, 14).Value = -7822

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3500.5625
$ws.Cells.Item(80, 9).Value = 3462.2307
$ws.Cells.Item(80, 10).Value = 3666.6667
$ws.Cells.Item(80, 11).Value = 3462.2307
$ws.Cells.Item(80, 12).Value = 3666.6667
$ws.Cells.Item(80, 13).Value = -2464.2307
$ws.Cells.Item(80, 14).Value = -5662.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 3500.5625
$ws.Cells.Item(83, 9).Value = 3462.2307
$ws.Cells.Item(83, 10).Value = 3666.6667
$ws.Cells.Item(83, 11).Value = 17311.1535
$ws.Cells.Item(83, 12).Value = 18333.3335
$ws.Cells.Item(83, 13).Value = -12319.1535
$ws.Cells.Item(83, 14).Value = -28317.3335

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2450.7585
$ws.Cells.Item(132, 9).Value = 1967.6316
$ws.Cells.Item(132, 10).Value = 3368.7
$ws.Cells.Item(132, 11).Value = 5902.8948
$ws.Cells.Item(132, 12).Value = 10106.1
$ws.Cells.Item(132, 13).Value = -3372.8948
$ws.Cells.Item(132, 14).Value = -15166.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(135, 8).Value = 54082.5
$ws.Cells.Item(135, 10).Value = 54082.5
$ws.Cells.Item(135, 12).Value = 54082.5
$ws.Cells.Item(135, 14).Value = -64222.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3626.9443
$ws.Cells.Item(7, 9).Value = 3603.125
$ws.Cells.Item(7, 10).Value = 3646
$ws.Cells.Item(7, 11).Value = 3603.125
$ws.Cells.Item(7, 12).Value = 3646
$ws.Cells.Item(7, 13).Value = -3491.125
$ws.Cells.Item(7, 14).Value = -3870

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 8259.666999999999
$ws.Cells.Item(40, 9).Value = 9389.5
$ws.Cells.Item(40, 11).Value = 9389.5
$ws.Cells.Item(40, 13).Value = -9253.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(126, 8).Value = 3626.9443
$ws.Cells.Item(126, 9).Value = 3603.125
$ws.Cells.Item(126, 10).Value = 3646
$ws.Cells.Item(126, 11).Value = 10809.375
$ws.Cells.Item(126, 12).Value = 10938
$ws.Cells.Item(126, 13).Value = -8339.375
$ws.Cells.Item(126, 14).Value = -15878

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1509.8462
$ws.Cells.Item(136, 9).Value = 1315.4783
$ws.Cells.Item(136, 10).Value = 3000
$ws.Cells.Item(136, 11).Value = 3946.4349
$ws.Cells.Item(136, 12).Value = 9000
$ws.Cells.Item(136, 13).Value = -1396.4349
$ws.Cells.Item(136, 14).Value = -14100
